$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F109").Value = 16
$ws.Range("G109").Value = 541.6
$ws.Range("F118").Value = 109
$ws.Range("G118").Value = 7664.88
$ws.Range("F120").Value = 197
$ws.Range("G120").Value = 27648.95
$ws.Range("F141").Value = 640
$ws.Range("G141").Value = 12454.4
$ws.Range("B143").Value = 369951.99
$ws.Range("F193").Value = 110
$ws.Range("G193").Value = 5442.8
$ws.Range("F199").Value = 315
$ws.Range("G199").Value = 6199.2
$ws.Range("B205").Value = 36372.58
$ws.Range("F215").Value = 28
$ws.Range("G215").Value = 1285.2
$ws.Range("B216").Value = 8199.07
$ws.Range("F251").Value = 8
$ws.Range("G251").Value = 895.04
$ws.Range("B257").Value = 6171.29
$ws.Range("F283").Value = 140
$ws.Range("G283").Value = 11876.2
$ws.Range("F284").Value = 57
$ws.Range("G284").Value = 4401.54
$ws.Range("B288").Value = 27392.3
$ws.Range("F322").Value = 26
$ws.Range("G322").Value = 525.46
$ws.Range("F328").Value = 16
$ws.Range("G328").Value = 969.28
$ws.Range("F341").Value = 36
$ws.Range("G341").Value = 727.5599999999999
$ws.Range("F349").Value = 19
$ws.Range("G349").Value = 807.12
$ws.Range("F363").Value = 65
$ws.Range("G363").Value = 2239.9
$ws.Range("B369").Value = 162351.4
$ws.Range("B374").Value = 57077
$ws.Range("D374").Value = 93.08
$ws.Range("E374").Value = 111.2
$ws.Range("F374").Value = 1
$ws.Range("G374").Value = 93.08
$ws.Range("B375").Value = 61610
$ws.Range("D375").Value = 102.71
$ws.Range("E375").Value = 122.71
$ws.Range("F375").Value = 398
$ws.Range("G375").Value = 40878.58
$ws.Range("F383").Value = 32
$ws.Range("G383").Value = 5077.44
$ws.Range("F389").Value = 8
$ws.Range("G389").Value = 207.52
$ws.Range("F417").Value = 70
$ws.Range("G417").Value = 7286.3
$ws.Range("F420").Value = 83
$ws.Range("G420").Value = 6624.23
$ws.Range("F429").Value = 286
$ws.Range("G429").Value = 16911.18
$ws.Range("F434").Value = 917
$ws.Range("G434").Value = 19889.73
$ws.Range("F435").Value = 564
$ws.Range("G435").Value = 3389.64
$ws.Range("F436").Value = 55
$ws.Range("G436").Value = 4570.5
$ws.Range("F437").Value = 35
$ws.Range("G437").Value = 18408.25
$ws.Range("F441").Value = 333
$ws.Range("G441").Value = 23413.23
$ws.Range("F447").Value = 71
$ws.Range("G447").Value = 4376.44
$ws.Range("F449").Value = 23
$ws.Range("G449").Value = 2169.82
$ws.Range("F450").Value = 135
$ws.Range("G450").Value = 2743.2
$ws.Range("F452").Value = 7
$ws.Range("G452").Value = 288.26
$ws.Range("B454").Value = 706700.47
$ws.Range("F458").Value = 5
$ws.Range("G458").Value = 808.65
$ws.Range("B469").Value = 71651.87
$ws.Range("F521").Value = 165
$ws.Range("G521").Value = 4907.1
$ws.Range("B529").Value = 80626.88
$ws.Range("F544").Value = 49
$ws.Range("G544").Value = 14387.38
$ws.Range("B558").Value = 89356.12
$ws.Range("F603").Value = 23
$ws.Range("G603").Value = 800.63
$ws.Range("F604").Value = 79
$ws.Range("G604").Value = 2749.99
$ws.Range("B607").Value = 15928.34
$ws.Range("F653").Value = 1
$ws.Range("G653").Value = 6117.08
$ws.Range("B654").Value = 11589.89
$ws.Range("F657").Value = 684
$ws.Range("G657").Value = 4685.4
$ws.Range("B664").Value = 51466.56
$ws.Range("F739").Value = 72
$ws.Range("G739").Value = 1972.8
$ws.Range("F743").Value = 44
$ws.Range("G743").Value = 2576.64
$ws.Range("F744").Value = 99
$ws.Range("G744").Value = 7281.45
$ws.Range("F745").Value = 108
$ws.Range("G745").Value = 13265.64
$ws.Range("B748").Value = 101834.18
$ws.Range("F752").Value = 173
$ws.Range("G752").Value = 22585.15
$ws.Range("B757").Value = 86300.53
$ws.Range("F781").Value = 46
$ws.Range("G781").Value = 1464.18
$ws.Range("B791").Value = 58828.12
$ws.Range("F820").Value = 71
$ws.Range("G820").Value = 6514.25
$ws.Range("B826").Value = 34402.91
$ws.Range("F850").Value = 46
$ws.Range("G850").Value = 4034.2
$ws.Range("B853").Value = 10045.42
$ws.Range("F885").Value = 217
$ws.Range("G885").Value = 31248
$ws.Range("F886").Value = 335
$ws.Range("G886").Value = 40437.85
$ws.Range("B888").Value = 152931.01
$ws.Range("F897").Value = 52
$ws.Range("G897").Value = 4430.4
$ws.Range("F915").Value = 31
$ws.Range("G915").Value = 986.11
$ws.Range("B920").Value = 83816.95
$ws.Range("F922").Value = 64
$ws.Range("G922").Value = 6885.12
$ws.Range("B929").Value = 37017.25
$ws.Range("F932").Value = 1
$ws.Range("G932").Value = 3650.99
$ws.Range("B939").Value = 224795.88
$ws.Range("F971").Value = 4
$ws.Range("G971").Value = 749.6
$ws.Range("B978").Value = 14781.67
$ws.Range("F983").Value = 117
$ws.Range("G983").Value = 17299.62
$ws.Range("F985").Value = 136
$ws.Range("G985").Value = 17486.88
$ws.Range("B986").Value = 79379.02
$ws.Range("B993").Value = 4319155.38
$ws.Range("B994").Value = 4319155.38
